$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102 - 20/06 update
$ws.Range("A102").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("A102").Value = 44367
$ws.Range("B102").Value = 734997
$ws.Range("C102").Value = 694462
$ws.Range("D102").Value = 270813
$ws.Range("E102").Value = 219138

# Row 103 - 21/06 update
$ws.Range("A103").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("A103").Value = 44368
$ws.Range("B103").Value = 734847
$ws.Range("C103").Value = 716309
$ws.Range("D103").Value = 270813
$ws.Range("E103").Value = 219298

[void]$ws.Range("B103").Select()
